# Junction_Flooding_423: refresh data rows, drop trailing row, widen columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace data rows 2-5 with refreshed readings ---
$row2 = New-Object 'object[,]' 1,34
$row2[0,0] = 45104.50694444445
$row2[0,1] = 11.794
$row2[0,2] = 7.946
$row2[0,3] = 3.447
$row2[0,4] = 25.671
$row2[0,5] = 18.994
$row2[0,6] = 8.994999999999999
$row2[0,7] = 26.449
$row2[0,8] = 14.543
$row2[0,9] = 5.683
$row2[0,10] = 8.23
$row2[0,11] = 10.12
$row2[0,12] = 10.985
$row2[0,13] = 3.013
$row2[0,14] = 9.398999999999999
$row2[0,15] = 12.777
$row2[0,16] = 8.539999999999999
$row2[0,17] = 2.71
$row2[0,18] = 1.486
$row2[0,19] = 135.674
$row2[0,20] = 25.95
$row2[0,21] = 8.676
$row2[0,22] = 16.478
$row2[0,23] = 8.538
$row2[0,24] = 2.569
$row2[0,25] = 14.652
$row2[0,26] = 7.663
$row2[0,27] = 7.143
$row2[0,28] = 8.202
$row2[0,29] = 10.671
$row2[0,30] = 2.707
$row2[0,31] = 24.031
$row2[0,32] = 4.513
$row2[0,33] = 10.846
$ws.Range("A2:AH2").Value = $row2

$row3 = New-Object 'object[,]' 1,34
$row3[0,0] = 45104.51388888889
$row3[0,1] = 14.249
$row3[0,2] = 10.338
$row3[0,3] = 1.687
$row3[0,4] = 31.276
$row3[0,5] = 24.654
$row3[0,6] = 11.071
$row3[0,7] = 41.835
$row3[0,8] = 17.452
$row3[0,9] = 7.48
$row3[0,10] = 10.846
$row3[0,11] = 12.5
$row3[0,12] = 13.419
$row3[0,13] = 3.623
$row3[0,14] = 11.279
$row3[0,15] = 15.823
$row3[0,16] = 9.869
$row3[0,17] = 1.342
$row3[0,18] = 0.929
$row3[0,19] = 164.347
$row3[0,20] = 31.527
$row3[0,21] = 10.411
$row3[0,22] = 20.764
$row3[0,23] = 10.85
$row3[0,24] = 2.023
$row3[0,25] = 21.046
$row3[0,26] = 9.196
$row3[0,27] = 8.337
$row3[0,28] = 9.750999999999999
$row3[0,29] = 13.118
$row3[0,30] = 1.141
$row3[0,31] = 38.438
$row3[0,32] = 5.675
$row3[0,33] = 13.016
$ws.Range("A3:AH3").Value = $row3

$row4 = New-Object 'object[,]' 1,34
$row4[0,0] = 45104.52083333334
$row4[0,1] = 12.84
$row4[0,2] = 9.439
$row4[0,3] = 1.184
$row4[0,4] = 28.189
$row4[0,5] = 22.479
$row4[0,6] = 10.016
$row4[0,7] = 40.03
$row4[0,8] = 15.707
$row4[0,9] = 6.831
$row4[0,10] = 9.936999999999999
$row4[0,11] = 11.289
$row4[0,12] = 12.094
$row4[0,13] = 3.261
$row4[0,14] = 10.151
$row4[0,15] = 14.319
$row4[0,16] = 8.789999999999999
$row4[0,17] = 0.92
$row4[0,18] = 0.7
$row4[0,19] = 147.187
$row4[0,20] = 28.414
$row4[0,21] = 9.369999999999999
$row4[0,22] = 18.841
$row4[0,23] = 9.859999999999999
$row4[0,24] = 1.689
$row4[0,25] = 19.538
$row4[0,26] = 8.276
$row4[0,27] = 7.453
$row4[0,28] = 8.734
$row4[0,29] = 11.852
$row4[0,30] = 0.724
$row4[0,31] = 36.627
$row4[0,32] = 5.156
$row4[0,33] = 11.715
$ws.Range("A4:AH4").Value = $row4

$row5 = New-Object 'object[,]' 1,34
$row5[0,0] = 45104.52777777778
$row5[0,1] = 8.539999999999999
$row5[0,2] = 6.28
$row5[0,3] = 0.83
$row5[0,4] = 18.8
$row5[0,5] = 14.91
$row5[0,6] = 6.66
$row5[0,7] = 28.29
$row5[0,8] = 10.47
$row5[0,9] = 4.55
$row5[0,10] = 6.56
$row5[0,11] = 7.53
$row5[0,12] = 8.09
$row5[0,13] = 2.18
$row5[0,14] = 6.77
$row5[0,15] = 9.539999999999999
$row5[0,16] = 5.9
$row5[0,17] = 0.68
$row5[0,18] = 0.48
$row5[0,19] = 95.69
$row5[0,20] = 19
$row5[0,21] = 6.25
$row5[0,22] = 12.57
$row5[0,23] = 6.56
$row5[0,24] = 1.16
$row5[0,25] = 13.56
$row5[0,26] = 5.52
$row5[0,27] = 4.99
$row5[0,28] = 5.84
$row5[0,29] = 7.91
$row5[0,30] = 0.53
$row5[0,31] = 25.94
$row5[0,32] = 3.42
$row5[0,33] = 7.81
$ws.Range("A5:AH5").Value = $row5

# --- Drop the old trailing row (was row 6; dataset now ends at row 5) ---
$ws.Rows.Item(6).Delete()

# --- Widen columns whose custom accuracy now needs an extra character ---
$ws.Columns.Item(2).ColumnWidth = 8 - (5/6)  # B: 7 -> 8
$ws.Columns.Item(3).ColumnWidth = 8 - (5/6)  # C: 7 -> 8
$ws.Columns.Item(5).ColumnWidth = 8 - (5/6)  # E: 7 -> 8
$ws.Columns.Item(6).ColumnWidth = 8 - (5/6)  # F: 7 -> 8
$ws.Columns.Item(7).ColumnWidth = 8 - (5/6)  # G: 7 -> 8
$ws.Columns.Item(8).ColumnWidth = 8 - (5/6)  # H: 7 -> 8
$ws.Columns.Item(9).ColumnWidth = 8 - (5/6)  # I: 7 -> 8
$ws.Columns.Item(11).ColumnWidth = 8 - (5/6)  # K: 7 -> 8
$ws.Columns.Item(12).ColumnWidth = 8 - (5/6)  # L: 7 -> 8
$ws.Columns.Item(13).ColumnWidth = 8 - (5/6)  # M: 7 -> 8
$ws.Columns.Item(15).ColumnWidth = 8 - (5/6)  # O: 7 -> 8
$ws.Columns.Item(16).ColumnWidth = 8 - (5/6)  # P: 7 -> 8
$ws.Columns.Item(20).ColumnWidth = 9 - (5/6)  # T: 8 -> 9
$ws.Columns.Item(21).ColumnWidth = 8 - (5/6)  # U: 7 -> 8
$ws.Columns.Item(22).ColumnWidth = 8 - (5/6)  # V: 7 -> 8
$ws.Columns.Item(23).ColumnWidth = 8 - (5/6)  # W: 7 -> 8
$ws.Columns.Item(26).ColumnWidth = 8 - (5/6)  # Z: 7 -> 8
$ws.Columns.Item(30).ColumnWidth = 8 - (5/6)  # AD: 7 -> 8
$ws.Columns.Item(32).ColumnWidth = 8 - (5/6)  # AF: 7 -> 8
$ws.Columns.Item(34).ColumnWidth = 8 - (5/6)  # AH: 7 -> 8
